# This script expands the last paragraph of the document ("Goal: Create a
# year column.") into a long section covering extraction of a time/year
# column from a date object and a new "count:" section with a ddply example,
# matching the target revision.
#
# Approach: Word's Range.InsertXML() accepts a full WordprocessingML package
# fragment and inserts it "as typed" at the target Range: every *complete*
# <w:p>...</w:p> paragraph in the fragment becomes its own new paragraph,
# while a trailing *unclosed* paragraph's runs are spliced into the
# paragraph that owns the target Range (preserving that paragraph's identity
# - including the "_GoBack" bookmark sitting at its end). So the fragment
# below is built as many whole new paragraphs followed by one more "open"
# paragraph (closed here just to keep the fragment well-formed XML) whose
# runs get merged in front of the pre-existing "Goal: " / "Create a year
# column. " runs and the bookmark.
#
# After the splice, the stale original "Goal: Create a year column. " text
# (now sitting after our new runs, still before the bookmark) is removed
# with a plain Find/Replace.

$d = $word.ActiveDocument
$lastPara = $d.Paragraphs.Last
$insertionPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Goal: Create a </w:t></w:r><w:r><w:t>time</w:t></w:r><w:r><w:t xml:space="preserve"> column. </w:t></w:r><w:r><w:t xml:space="preserve">The two challenges that may be associated with this are if the date is provided as a date formatted object (example 01/01/2015) or if there are multiple samples per year. If the latter is the issue, the appropriate </w:t></w:r><w:r><w:t>temporal</w:t></w:r><w:r><w:t xml:space="preserve"> scale must be determined prior to formatting the dataset (see Allen, Ethan, or myself).</w:t></w:r><w:r><w:t xml:space="preserve"> Once the time scale is determined, data are reported as a decimal year and simply requires a bit of math (for example, if sampling was done monthly and a sample was taken on 1 Mar 2015, the time of the sample would be 2015 + 3/12 as March is the third month of the year).</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:rPr><w:b/><w:i/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:i/></w:rPr><w:t>Extracting year from a date object:</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">Convert the date column to an R formatted date (in this case pretending that our unformatted dataset contains a column called </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>record_date</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>):</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:t>date</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> = </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>strptime</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>example_df</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> $</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>record_date</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, '%m/ %d/ %y')</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Add a sampling year line (summarize by year):</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>example</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>_df</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> $year = </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>as.numeric</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(format(date, '%Y'))</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:rPr><w:b/><w:sz w:val="32"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/><w:sz w:val="32"/></w:rPr><w:t>count</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b/><w:sz w:val="32"/></w:rPr><w:t>:</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Goal: Summarize the dataset to the count of individuals per species, site, and year for a given dataset. To do so, we will use Hadley Wickham’s “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ddply</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">” function in the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>plyr</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> package.</w:t></w:r><w:r><w:t xml:space="preserve"> Below is an example in which there is a count column that must be summarized.</w:t></w:r></w:p><w:p/><w:p><w:proofErr w:type="gramStart"/><w:r><w:t>example</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">_df2 = </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ddply</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>example_df</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, .(site, year, species</w:t></w:r><w:r><w:t xml:space="preserve">), summarize, count = </w:t></w:r><w:r><w:t>sum(count</w:t></w:r><w:r><w:t>))</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$insertionPoint.InsertXML($xml)

$null = $d.Content.Find.Execute("Goal: Create a year column. ", $false, $false, $false, $false, $false, $true, 1, $false, "", 2)
